$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13. This shifts the former rows 13-21
# down to 14-22 (carrying their row heights with them), and creates a
# new blank row 13 that will hold the "Docentes responsaveis" data.
$ws.Rows.Item(13).Insert()

# Update the cells whose text changed (new/rearranged content).
# Row 10
$ws.Range("B10").Value = "Introduzir e estabelecer aos alunos os princípios básicos a serem utilizados em todos os processos que envolvam a conservação de massa e energia. Esta disciplina propicia a realização de balanços globais de massa e energia em diferentes processos químicos evidenciando a importância da aplicação desta metodologia no projeto e otimização de processos químicos industriais."
$ws.Range("C10").Value = "Introduzir e estabelecer aos alunos os princípios básicos a serem utilizados em todos os processos que envolvam a conservação de massa e energia. Esta disciplina propicia a realização de balanços globais de massa e energia em diferentes processos químicos evidenciando a importância da aplicação desta metodologia no projeto e otimização de processos químicos industriais."

# Row 13
$ws.Range("B13").Value = "5817045 - Elisângela de Jesus Cândido Moraes"
$ws.Range("C13").Value = "5817045 - Elisângela de Jesus Cândido Moraes"

# Row 14
$ws.Range("B14").Value = "Introdução aos cálculos em Engenharia Química; Balanços materiais; Balanços de energia; Balanços material e energético combinados; Balanços em processos no estado transiente."
$ws.Range("C14").Value = "Introdução aos cálculos em Engenharia Química; Balanços materiais; Balanços de energia; Balanços material e energético combinados; Balanços em processos no estado transiente."

# Row 16
$ws.Range("B16").Value = "1 - Balanços Materiais- Introdução aos Balanços Materiais- Balanços Materiais que não envolvem reações químicas.- Balanços Materiais envolvendo reações químicas- Balanços Materiais com recirculação (reciclo e Bypass).2 - Balanços de Energia- Definições e conceitos. Formas de energia, calor, entalpia, valores de entalpia e capacidade calorífica.- Balanços de energia que não envolvem reações químicas.- Balanços de energia envolvendo reações químicas.3 - Balanços de Massa e de Energia Combinados- Aplicação combinada dos balanços de massa e energia em processos tais como umidificação, dissolução, processos de mistura, etc.4 - Balanços de Massa e de Energia em Processos de Regime Transiente"
$ws.Range("C16").Value = "1 - Balanços Materiais- Introdução aos Balanços Materiais- Balanços Materiais que não envolvem reações químicas.- Balanços Materiais envolvendo reações químicas- Balanços Materiais com recirculação (reciclo e Bypass).2 - Balanços de Energia- Definições e conceitos. Formas de energia, calor, entalpia, valores de entalpia e capacidade calorífica.- Balanços de energia que não envolvem reações químicas.- Balanços de energia envolvendo reações químicas.3 - Balanços de Massa e de Energia Combinados- Aplicação combinada dos balanços de massa e energia em processos tais como umidificação, dissolução, processos de mistura, etc.4 - Balanços de Massa e de Energia em Processos de Regime Transiente"

# Row 19
$ws.Range("B19").Value = "Provas escritas; -participação e conteúdo de trabalho e seminário;"
$ws.Range("C19").Value = "Provas escritas; -participação e conteúdo de trabalho e seminário;"

# Row 20
$ws.Range("B20").Value = "Média Final = (Prova1 + 2xProva2 + Nota de Trabalho) / 4`nMédia final mínima de aprovação = 5,0"
$ws.Range("C20").Value = "Média Final = (Prova1 + 2xProva2 + Nota de Trabalho) / 4`nMédia final mínima de aprovação = 5,0"

# Row 21
$ws.Range("B21").Value = "(Prova escrita + Média Final)/2`nNota Final mínima para aprovação= 5,0"
$ws.Range("C21").Value = "(Prova escrita + Média Final)/2`nNota Final mínima para aprovação= 5,0"

# Row 22
$ws.Range("B22").Value = "HIMMELBLAU, David M.  Eng. Química princípios e cálculos.  7. ed. LTC Editora,2006.`nGOMIDE, R.  Estequiometria Industrial. 3.ed. São Paulo: Ed. do Autor, 1984. `nFELDER, R.M; ROUSSEAU, R.W. Princípios elementares dos processos químicos. 3. ed. Rio de Janeiro: LTC Editora, 2005`nMOUYEN, O.A.; WATSON, K. M. AND RAGATZ, R.A.  Princípios dos processos químicos.  Lisboa: Lopes da Silva Editora, 2005. v. 1 `nCREMASCO, M. A. Fundamentos de transferência de massa. 1.ed. Campinas: Editora da UNICAMP, 1998.`nBRASIL, N. I. Introdução à Engenharia Química.  2. ed. Rio de Janeiro: Editora Interciência , 2004."
$ws.Range("C22").Value = "HIMMELBLAU, David M.  Eng. Química princípios e cálculos.  7. ed. LTC Editora,2006.`nGOMIDE, R.  Estequiometria Industrial. 3.ed. São Paulo: Ed. do Autor, 1984. `nFELDER, R.M; ROUSSEAU, R.W. Princípios elementares dos processos químicos. 3. ed. Rio de Janeiro: LTC Editora, 2005`nMOUYEN, O.A.; WATSON, K. M. AND RAGATZ, R.A.  Princípios dos processos químicos.  Lisboa: Lopes da Silva Editora, 2005. v. 1 `nCREMASCO, M. A. Fundamentos de transferência de massa. 1.ed. Campinas: Editora da UNICAMP, 1998.`nBRASIL, N. I. Introdução à Engenharia Química.  2. ed. Rio de Janeiro: Editora Interciência , 2004."
